# Refresh cryptocurrency price/volume snapshot (GitHub Actions scheduled update).
# Prices in column D are plain text (not numbers) in the source data -- some of
# them look like ordinary decimals (e.g. "536.77") which Excel would otherwise
# auto-convert to a Number on assignment. Set-TextValue forces those through as
# literal text (matching the original inlineStr cells) and then restores the
# cell's default (unstyled) formatting so no stray style gets left behind.
function Set-TextValue($sheet, $addr, $val) {
    $c = $sheet.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

Set-TextValue $ws "D2" "60.027.40"
$ws.Range("E2").Value = "  +2.23%  "
Set-TextValue $ws "D3" "3.194.16"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("E4").Value = "  +0.00%  "
Set-TextValue $ws "D5" "536.77"
$ws.Range("E5").Value = "  +1.01%  "
Set-TextValue $ws "D6" "145.45"
$ws.Range("E6").Value = "  +4.09%  "
$ws.Range("E7").Value = "  +0.01%  "
Set-TextValue $ws "D8" "0.531"
$ws.Range("E8").Value = "  -0.37%  "
$ws.Range("E9").Value = "  -0.08%  "
Set-TextValue $ws "D10" "0.113"
$ws.Range("E10").Value = "  +2.07%  "
$ws.Range("E11").Value = "  -0.74%  "
Set-TextValue $ws "D12" "3.745.56"
$ws.Range("E12").Value = "  +1.29%  "
$ws.Range("E13").Value = "  -2.00%  "
Set-TextValue $ws "D14" "25.85"
$ws.Range("E14").Value = "  -0.71%  "
Set-TextValue $ws "D15" "0.0000173"
$ws.Range("E15").Value = "  +0.12%  "
Set-TextValue $ws "D16" "60.046.77"
$ws.Range("E16").Value = "  +2.18%  "
Set-TextValue $ws "D17" "3.192.98"
$ws.Range("E17").Value = "  +1.21%  "
Set-TextValue $ws "D18" "6.28"
$ws.Range("E18").Value = "  +0.21%  "
Set-TextValue $ws "D19" "13.31"
$ws.Range("E19").Value = "  +1.59%  "
$ws.Range("E20").Value = "  +0.51%  "
Set-TextValue $ws "D21" "369.70"
$ws.Range("E21").Value = "  -0.98%  "
$ws.Range("E22").Value = "  +0.39%  "
Set-TextValue $ws "D23" "0.523"
$ws.Range("E23").Value = "  +0.13%  "
Set-TextValue $ws "D24" "69.47"
$ws.Range("E24").Value = "  -0.91%  "
$ws.Range("E25").Value = "  +0.95%  "
Set-TextValue $ws "D26" "8.66"
$ws.Range("E26").Value = "  +5.04%  "
Set-TextValue $ws "D27" "0.999"
$ws.Range("E27").Value = "  -0.18%  "
Set-TextValue $ws "D28" "0.0₃0876"
$ws.Range("E28").Value = "  +1.08%  "
Set-TextValue $ws "D29" "22.49"
$ws.Range("E29").Value = "  +1.53%  "
$ws.Range("E30").Value = "  +0.81%  "
Set-TextValue $ws "D31" "6.13"
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("E32").Value = "  +1.57%  "
$ws.Range("E33").Value = "  +1.72%  "
Set-TextValue $ws "D34" "6.58"
$ws.Range("E34").Value = "  +4.95%  "
Set-TextValue $ws "D35" "156.48"
$ws.Range("E35").Value = "  -1.72%  "
Set-TextValue $ws "D36" "1.36"
$ws.Range("E36").Value = "  +1.49%  "
Set-TextValue $ws "D37" "2.845.13"
$ws.Range("E37").Value = "  +7.55%  "
Set-TextValue $ws "D38" "26.33"
$ws.Range("E38").Value = "  +4.58%  "
Set-TextValue $ws "D39" "0.0706"
$ws.Range("E39").Value = "  +2.94%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws "D40" "0.0309"
$ws.Range("E40").Value = "  +8.18%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws "D41" "1.67"
$ws.Range("E41").Value = "  -0.61%  "
Set-TextValue $ws "D42" "4.24"
$ws.Range("E42").Value = "  -0.32%  "
Set-TextValue $ws "D43" "39.93"
$ws.Range("E43").Value = "  +2.73%  "
Set-TextValue $ws "D44" "0.719"
$ws.Range("E44").Value = "  +1.42%  "
$ws.Range("E45").Value = "  +1.36%  "
Set-TextValue $ws "D46" "3.237.33"
$ws.Range("E46").Value = "  +1.27%  "
Set-TextValue $ws "D47" "0.986"
$ws.Range("E47").Value = "  +0.25%  "
$ws.Range("E48").Value = "  -1.10%  "
Set-TextValue $ws "D49" "20.73"
$ws.Range("E49").Value = "  +1.55%  "
$ws.Range("E50").Value = "  +4.64%  "
$ws.Range("E51").Value = "  +0.01%  "
